# Insert a new data row at row 83, shifting existing rows 83..213 down to 84..214.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("83:83").Insert()

$ws.Range("A83").Value = 8
$ws.Range("B83").Value = "Terminal La Palmera de La Serena"
$ws.Range("C83").Value = "Coquimbo"
$ws.Range("D83").Value = 44495
$ws.Range("D83").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E83").Value = 4
$ws.Range("F83").Value = 100114013
$ws.Range("G83").Value = "Zanahoria"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 660
$ws.Range("K83").Value = 6500
$ws.Range("L83").Value = 7000
$ws.Range("M83").Value = 6750
$ws.Range("N83").Value = "$/saco 20 kilos"
$ws.Range("O83").Value = "Provincia del Elquí"
$ws.Range("P83").Value = 338
$ws.Range("Q83").Value = 20
$ws.Range("R83").Value = "Hortaliza"
